$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (osimertinib, gmax=5)
$ws.Range("B2").Value = 0.2
$ws.Range("C2").Value = -0.1
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1

# Row 3 (osimertinib, gmax=5)
$ws.Range("B3").Value = 0.15
$ws.Range("C3").Value = -0.1
$ws.Range("D3").Value = 10
$ws.Range("E3").Value = 1

# Row 4 (osimertinib, gmax=5)
$ws.Range("B4").Value = 0.2
$ws.Range("C4").Value = -0.1
$ws.Range("D4").Value = 100
$ws.Range("E4").Value = 1

# Row 5 (osimertinib, gmax=5)
$ws.Range("B5").Value = 0.1
$ws.Range("C5").Value = -0.1
$ws.Range("D5").Value = 1000
$ws.Range("E5").Value = 1

# Row 6 (gefitinib, gmax=7)
$ws.Range("B6").Value = 0.2
$ws.Range("C6").Value = -0.1
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1

# Row 7 (gefitinib, gmax=7)
$ws.Range("B7").Value = 0.15
$ws.Range("C7").Value = -0.1
$ws.Range("D7").Value = 100
$ws.Range("E7").Value = 1

# Row 8 (gefitinib, gmax=7)
$ws.Range("B8").Value = 0.2
$ws.Range("C8").Value = -0.1
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 1

# Row 9 (gefitinib, gmax=7)
$ws.Range("B9").Value = 0.1
$ws.Range("C9").Value = -0.1
$ws.Range("D9").Value = 10
$ws.Range("E9").Value = 1

# Update view/selection to match the saved workbook state (no scrolled topLeftCell,
# active cell D8 selected)
$ws.Range("D8").Select()
